$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 17.41485066666667
$ws.Range("H2").Value = 52.24455200000001
$ws.Range("I2").Value = 0.1047285618770465
$ws.Range("J2").Value = 0.1047285618770465
$ws.Range("M2").Value = 0.1419263333333333
$ws.Range("N2").Value = 0.425779
$ws.Range("O2").Value = 0.002583058778296354
$ws.Range("P2").Value = 0.002583058778296354
$ws.Range("Q2").Value = 2.471625900667556
$ws.Range("R2").Value = 22.244633106008
$ws.Range("S2").Value = 0.0002705200310948579
$ws.Range("T2").Value = 0.0002705200310948579
# Row 3
$ws.Range("G3").Value = 17.41485066666667
$ws.Range("H3").Value = 52.24455200000001
$ws.Range("I3").Value = 0.1047285618770465
$ws.Range("J3").Value = 0.1047285618770465
$ws.Range("O3").Value = 0.001399682868699959
$ws.Range("P3").Value = 0.001399682868699959
$ws.Range("Q3").Value = 1.339300700420445
$ws.Range("R3").Value = 12.053706303784
$ws.Range("S3").Value = 0.0001465867739228856
$ws.Range("T3").Value = 0.0001465867739228856
# Row 4
$ws.Range("G4").Value = 17.41485066666667
$ws.Range("H4").Value = 52.24455200000001
$ws.Range("I4").Value = 0.1047285618770465
$ws.Range("J4").Value = 0.1047285618770465
$ws.Range("M4").Value = 2.613991
$ws.Range("N4").Value = 7.841973
$ws.Range("O4").Value = 0.04757462720522382
$ws.Range("P4").Value = 0.04757462720522382
$ws.Range("Q4").Value = 45.52226290901067
$ws.Range("R4").Value = 409.7003661810961
$ws.Range("S4").Value = 0.004982422289039704
$ws.Range("T4").Value = 0.004982422289039704
# Row 5
$ws.Range("G5").Value = 17.41485066666667
$ws.Range("H5").Value = 52.24455200000001
$ws.Range("I5").Value = 0.1047285618770465
$ws.Range("J5").Value = 0.1047285618770465
$ws.Range("M5").Value = 52.11224233333333
$ws.Range("N5").Value = 156.336727
$ws.Range("O5").Value = 0.9484426311477799
$ws.Range("P5").Value = 0.9484426311477798
$ws.Range("Q5").Value = 907.5269181401451
$ws.Range("R5").Value = 8167.742263261304
$ws.Range("S5").Value = 0.09932903278298909
$ws.Range("T5").Value = 0.09932903278298907
# Row 6
$ws.Range("I6").Value = 0.1785014126970782
$ws.Range("J6").Value = 0.1785014126970782
$ws.Range("M6").Value = 0.1419263333333333
$ws.Range("N6").Value = 0.425779
$ws.Range("O6").Value = 0.002583058778296354
$ws.Range("P6").Value = 0.002583058778296354
$ws.Range("Q6").Value = 4.212687609000223
$ws.Range("R6").Value = 37.914188481002
$ws.Range("S6").Value = 0.000461079641005488
$ws.Range("T6").Value = 0.000461079641005488
# Row 7
$ws.Range("I7").Value = 0.1785014126970782
$ws.Range("J7").Value = 0.1785014126970782
$ws.Range("O7").Value = 0.001399682868699959
$ws.Range("P7").Value = 0.001399682868699959
$ws.Range("S7").Value = 0.0002498453693908416
$ws.Range("T7").Value = 0.0002498453693908417
# Row 8
$ws.Range("I8").Value = 0.1785014126970782
$ws.Range("J8").Value = 0.1785014126970782
$ws.Range("M8").Value = 2.613991
$ws.Range("N8").Value = 7.841973
$ws.Range("O8").Value = 0.04757462720522382
$ws.Range("P8").Value = 0.04757462720522382
$ws.Range("Q8").Value = 77.58903677075267
$ws.Range("R8").Value = 698.3013309367741
$ws.Range("S8").Value = 0.008492138164669299
$ws.Range("T8").Value = 0.008492138164669301
# Row 9
$ws.Range("I9").Value = 0.1785014126970782
$ws.Range("J9").Value = 0.1785014126970782
$ws.Range("M9").Value = 52.11224233333333
$ws.Range("N9").Value = 156.336727
$ws.Range("O9").Value = 0.9484426311477799
$ws.Range("P9").Value = 0.9484426311477798
$ws.Range("Q9").Value = 1546.80665947487
$ws.Range("R9").Value = 13921.25993527383
$ws.Range("S9").Value = 0.1692983495220126
$ws.Range("T9").Value = 0.1692983495220126
# Row 10
$ws.Range("G10").Value = 84.03051233333333
$ws.Range("H10").Value = 252.091537
$ws.Range("I10").Value = 0.5053385113032314
$ws.Range("J10").Value = 0.5053385113032314
$ws.Range("M10").Value = 0.1419263333333333
$ws.Range("N10").Value = 0.425779
$ws.Range("O10").Value = 0.002583058778296354
$ws.Range("P10").Value = 0.002583058778296354
$ws.Range("Q10").Value = 11.92614250359145
$ws.Range("R10").Value = 107.335282532323
$ws.Range("S10").Value = 0.001305319077633023
$ws.Range("T10").Value = 0.001305319077633023
# Row 11
$ws.Range("G11").Value = 84.03051233333333
$ws.Range("H11").Value = 252.091537
$ws.Range("I11").Value = 0.5053385113032314
$ws.Range("J11").Value = 0.5053385113032314
$ws.Range("O11").Value = 0.001399682868699959
$ws.Range("P11").Value = 0.001399682868699959
$ws.Range("Q11").Value = 6.462422571336555
$ws.Range("R11").Value = 58.161803142029
$ws.Range("S11").Value = 0.0007073136571654735
$ws.Range("T11").Value = 0.0007073136571654735
# Row 12
$ws.Range("G12").Value = 84.03051233333333
$ws.Range("H12").Value = 252.091537
$ws.Range("I12").Value = 0.5053385113032314
$ws.Range("J12").Value = 0.5053385113032314
$ws.Range("M12").Value = 2.613991
$ws.Range("N12").Value = 7.841973
$ws.Range("O12").Value = 0.04757462720522382
$ws.Range("P12").Value = 0.04757462720522382
$ws.Range("Q12").Value = 219.6550029647223
$ws.Range("R12").Value = 1976.895026682501
$ws.Range("S12").Value = 0.02404129128769401
$ws.Range("T12").Value = 0.02404129128769401
# Row 13
$ws.Range("G13").Value = 84.03051233333333
$ws.Range("H13").Value = 252.091537
$ws.Range("I13").Value = 0.5053385113032314
$ws.Range("J13").Value = 0.5053385113032314
$ws.Range("M13").Value = 52.11224233333333
$ws.Range("N13").Value = 156.336727
$ws.Range("O13").Value = 0.9484426311477799
$ws.Range("P13").Value = 0.9484426311477798
$ws.Range("Q13").Value = 4379.018422108822
$ws.Range("R13").Value = 39411.1657989794
$ws.Range("S13").Value = 0.4792845872807389
$ws.Range("T13").Value = 0.4792845872807389
# Row 14
$ws.Range("G14").Value = 35.158014
$ws.Range("H14").Value = 105.474042
$ws.Range("I14").Value = 0.2114315141226439
$ws.Range("J14").Value = 0.2114315141226439
$ws.Range("M14").Value = 0.1419263333333333
$ws.Range("N14").Value = 0.425779
$ws.Range("O14").Value = 0.002583058778296354
$ws.Range("P14").Value = 0.002583058778296354
$ws.Range("Q14").Value = 4.989848014302001
$ws.Range("R14").Value = 44.908632128718
$ws.Range("S14").Value = 0.0005461400285629848
$ws.Range("T14").Value = 0.0005461400285629847
# Row 15
$ws.Range("G15").Value = 35.158014
$ws.Range("H15").Value = 105.474042
$ws.Range("I15").Value = 0.2114315141226439
$ws.Range("J15").Value = 0.2114315141226439
$ws.Range("O15").Value = 0.001399682868699959
$ws.Range("P15").Value = 0.001399682868699959
$ws.Range("Q15").Value = 2.703850505346
$ws.Range("R15").Value = 24.334654548114
$ws.Range("S15").Value = 0.0002959370682207581
$ws.Range("T15").Value = 0.0002959370682207581
# Row 16
$ws.Range("G16").Value = 35.158014
$ws.Range("H16").Value = 105.474042
$ws.Range("I16").Value = 0.2114315141226439
$ws.Range("J16").Value = 0.2114315141226439
$ws.Range("M16").Value = 2.613991
$ws.Range("N16").Value = 7.841973
$ws.Range("O16").Value = 0.04757462720522382
$ws.Range("P16").Value = 0.04757462720522382
$ws.Range("Q16").Value = 91.902732173874
$ws.Range("R16").Value = 827.124589564866
$ws.Range("S16").Value = 0.0100587754638208
$ws.Range("T16").Value = 0.0100587754638208
# Row 17
$ws.Range("G17").Value = 35.158014
$ws.Range("H17").Value = 105.474042
$ws.Range("I17").Value = 0.2114315141226439
$ws.Range("J17").Value = 0.2114315141226439
$ws.Range("M17").Value = 52.11224233333333
$ws.Range("N17").Value = 156.336727
$ws.Range("O17").Value = 0.9484426311477799
$ws.Range("P17").Value = 0.9484426311477798
$ws.Range("Q17").Value = 1832.162945526726
$ws.Range("R17").Value = 16489.46650974053
$ws.Range("S17").Value = 0.2005306615620394
$ws.Range("T17").Value = 0.2005306615620393
